$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.00"
$ws.Range("E2").Value = "'-0.34%"
$ws.Range("D3").Value = "'39.77"
$ws.Range("E3").Value = "'1.99%"
$ws.Range("D4").Value = "'5.159"
$ws.Range("E4").Value = "'1.51%"
$ws.Range("D5").Value = "'0.08104"
$ws.Range("E5").Value = "'-1.07%"
$ws.Range("D6").Value = "'1.943"
$ws.Range("E6").Value = "'-4.21%"
$ws.Range("D7").Value = "'8.147"
$ws.Range("E7").Value = "'3.04%"
$ws.Range("E8").Value = "'-0.21%"
$ws.Range("D9").Value = "'0.1429"
$ws.Range("E9").Value = "'0.18%"
$ws.Range("D10").Value = "'0.1927"
$ws.Range("E10").Value = "'-0.97%"
$ws.Range("D11").Value = "'0.09160"
$ws.Range("E11").Value = "'-1.25%"
$ws.Range("D12").Value = "'0.03508"
$ws.Range("E12").Value = "'1.15%"
$ws.Range("D13").Value = "'0.09820"
$ws.Range("E13").Value = "'-0.25%"
$ws.Range("D14").Value = "'0.001403"
$ws.Range("E14").Value = "'-0.73%"
$ws.Range("D15").Value = "'0.005882"
$ws.Range("E15").Value = "'0.53%"
$ws.Range("E16").Value = "'2.34%"
$ws.Range("D17").Value = "'4.236"
$ws.Range("E17").Value = "'1.40%"
$ws.Range("D18").Value = "'3.356"
$ws.Range("E18").Value = "'-2.35%"
$ws.Range("D19").Value = "'0.3428"
$ws.Range("D20").Value = "'0.1321"
$ws.Range("E20").Value = "'1.45%"
$ws.Range("D21").Value = "'4.656"
$ws.Range("E21").Value = "'-3.81%"
$ws.Range("E22").Value = "'2.98%"
$ws.Range("D23").Value = "'0.04368"
$ws.Range("E23").Value = "'-2.49%"
$ws.Range("D24").Value = "'0.001230"
$ws.Range("E24").Value = "'-1.01%"
$ws.Range("D25").Value = "'0.004359"
$ws.Range("E25").Value = "'4.44%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'0.12%"
$ws.Range("D27").Value = "'0.0004003"
$ws.Range("E27").Value = "'-9.99%"
$ws.Range("D39").Value = "'0.02038"
$ws.Range("E39").Value = "'-3.66%"
$ws.Range("D40").Value = "'0.05067"
$ws.Range("E40").Value = "'-2.15%"
$ws.Range("D41").Value = "'0.007396"
$ws.Range("E41").Value = "'-0.96%"
$ws.Range("D42").Value = "'0.009770"
$ws.Range("E42").Value = "'-3.61%"
$ws.Range("E43").Value = "'-0.30%"
$ws.Range("D44").Value = "'0.002133"
$ws.Range("E44").Value = "'0.12%"
$ws.Range("D45").Value = "'0.009488"
$ws.Range("E45").Value = "'-1.93%"
$ws.Range("D46").Value = "'0.00006346"
$ws.Range("E46").Value = "'0.76%"
$ws.Range("E47").Value = "'0.05%"
$ws.Range("D48").Value = "'0.002731"
$ws.Range("E49").Value = "'-18.75%"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("E51").Value = "'0.05%"
